$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.023464926401103
$ws.Range("D2").Value = 1.029410939858912
$ws.Range("E2").Value = 1.024081151752958
$ws.Range("F2").Value = 1.036791413941509
$ws.Range("I2").Value = 1.032957945433059
$ws.Range("J2").Value = 1.02864499537528
$ws.Range("K2").Value = 1.0322249298917
$ws.Range("L2").Value = 1.02691070080248
$ws.Range("M2").Value = 1.03958414543502
$ws.Range("N2").Value = 1.013605090547963

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024223545733136
$ws.Range("D3").Value = 1.02999149016094
$ws.Range("E3").Value = 1.02472010550108
$ws.Range("F3").Value = 1.038922060843811
$ws.Range("I3").Value = 1.033162750708813
$ws.Range("J3").Value = 1.02904309398478
$ws.Range("K3").Value = 1.03261423283312
$ws.Range("L3").Value = 1.027357127712651
$ws.Range("M3").Value = 1.041520969180735
$ws.Range("N3").Value = 1.013737860397113

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024714854512626
$ws.Range("D4").Value = 1.030367385749939
$ws.Range("E4").Value = 1.025134314628816
$ws.Range("F4").Value = 1.04029575081282
$ws.Range("I4").Value = 1.033293984101119
$ws.Range("J4").Value = 1.029300444973271
$ws.Range("K4").Value = 1.032865678706261
$ws.Range("L4").Value = 1.027646071699587
$ws.Range("M4").Value = 1.042768900318437
$ws.Range("N4").Value = 1.013823661988932

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024921503107513
$ws.Range("D5").Value = 1.030525469233882
$ws.Range("E5").Value = 1.025308629793566
$ws.Range("F5").Value = 1.040872086696363
$ws.Range("I5").Value = 1.033348846787163
$ws.Range("J5").Value = 1.029408576233781
$ws.Range("K5").Value = 1.032971276493691
$ws.Range("L5").Value = 1.027767561335126
$ws.Range("M5").Value = 1.043292283191982
$ws.Range("N5").Value = 1.013859706640468

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024956206293864
$ws.Range("D6").Value = 1.030552015468326
$ws.Range("E6").Value = 1.025337908696501
$ws.Range("F6").Value = 1.040968788644349
$ws.Range("I6").Value = 1.033358040445941
$ws.Range("J6").Value = 1.029426728472747
$ws.Range("K6").Value = 1.032989000368796
$ws.Range("L6").Value = 1.027787960993099
$ws.Range("M6").Value = 1.043380089111712
$ws.Range("N6").Value = 1.013865757147089

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0247176153601
$ws.Range("D7").Value = 1.030369497846168
$ws.Range("E7").Value = 1.02513664312424
$ws.Range("F7").Value = 1.040303456376967
$ws.Range("I7").Value = 1.033294718387587
$ws.Range("J7").Value = 1.029301890061125
$ws.Range("K7").Value = 1.032867090141724
$ws.Range("L7").Value = 1.027647694981727
$ws.Range("M7").Value = 1.042775898654455
$ws.Range("N7").Value = 1.01382414372302

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023721215518605
$ws.Range("D8").Value = 1.029607089323332
$ws.Range("E8").Value = 1.02429693045878
$ws.Range("F8").Value = 1.037512525017196
$ws.Range("I8").Value = 1.033027427803917
$ws.Range("J8").Value = 1.0287795852839
$ws.Range("K8").Value = 1.032356591641543
$ws.Range("L8").Value = 1.027061556794217
$ws.Range("M8").Value = 1.040239822698363
$ws.Range("N8").Value = 1.013649983297792

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.021968773194031
$ws.Range("D9").Value = 1.028265499187254
$ws.Range("E9").Value = 1.022823143048933
$ws.Range("F9").Value = 1.0325550815297
$ws.Range("I9").Value = 1.03254651457306
$ws.Range("J9").Value = 1.027857351767498
$ws.Range("K9").Value = 1.031453516069558
$ws.Range("L9").Value = 1.026029313213571
$ws.Range("M9").Value = 1.035728938178789
$ws.Range("N9").Value = 1.013342259056363

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020802773585121
$ws.Range("D10").Value = 1.027372397909433
$ws.Range("E10").Value = 1.021844642702007
$ws.Range("F10").Value = 1.029221728645622
$ws.Range("I10").Value = 1.032219184766467
$ws.Range("J10").Value = 1.027241291507634
$ws.Range("K10").Value = 1.030849107775085
$ws.Range("L10").Value = 1.02534159096946
$ws.Range("M10").Value = 1.032691721289522
$ws.Range("N10").Value = 1.01313655849837

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020298436759206
$ws.Range("D11").Value = 1.026985987975731
$ws.Range("E11").Value = 1.021421908347578
$ws.Range("F11").Value = 1.027771199001914
$ws.Range("I11").Value = 1.032075840068023
$ws.Range("J11").Value = 1.026974239753384
$ws.Range("K11").Value = 1.030586832825034
$ws.Range("L11").Value = 1.025043910212445
$ws.Range("M11").Value = 1.031369078579223
$ws.Range("N11").Value = 1.013047358603807

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020111186731208
$ws.Range("D12").Value = 1.026842504944248
$ws.Range("E12").Value = 1.02126503159476
$ws.Range("F12").Value = 1.027231294917613
$ws.Range("I12").Value = 1.032022352565179
$ws.Range("J12").Value = 1.026875000881505
$ws.Range("K12").Value = 1.030489327659391
$ws.Range("L12").Value = 1.02493335491835
$ws.Range("M12").Value = 1.030876628541469
$ws.Range("N12").Value = 1.013014206322033

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020151348725443
$ws.Range("D13").Value = 1.026873280416193
$ws.Range("E13").Value = 1.021298675586404
$ws.Range("F13").Value = 1.027347157120056
$ws.Range("I13").Value = 1.032033836827507
$ws.Range("J13").Value = 1.026896289943493
$ws.Range("K13").Value = 1.030510246677512
$ws.Range("L13").Value = 1.024957068646476
$ws.Range("M13").Value = 1.030982313859854
$ws.Range("N13").Value = 1.013021318478653

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020282956913853
$ws.Range("D14").Value = 1.026974126666059
$ws.Range("E14").Value = 1.021408937890378
$ws.Range("F14").Value = 1.027726593274055
$ws.Range("I14").Value = 1.032071423736516
$ws.Range("J14").Value = 1.026966037533256
$ws.Range("K14").Value = 1.030578774743597
$ws.Range("L14").Value = 1.025034771336375
$ws.Range("M14").Value = 1.031328396399597
$ws.Range("N14").Value = 1.013044618622521

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020364056097423
$ws.Range("D15").Value = 1.027036267585675
$ws.Range("E15").Value = 1.021476893458024
$ws.Range("F15").Value = 1.027960227862675
$ws.Range("I15").Value = 1.032094550048536
$ws.Range("J15").Value = 1.027009005504271
$ws.Range("K15").Value = 1.030620985936498
$ws.Range("L15").Value = 1.025082648739738
$ws.Range("M15").Value = 1.031541474327802
$ws.Range("N15").Value = 1.013058972032887

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020836256379014
$ws.Range("D16").Value = 1.027398049207605
$ws.Range("E16").Value = 1.021872718555943
$ws.Range("F16").Value = 1.029317841468261
$ws.Range("I16").Value = 1.032228664085687
$ws.Range("J16").Value = 1.02725900868057
$ws.Range("K16").Value = 1.030866502246523
$ws.Range("L16").Value = 1.025361349364167
$ws.Range("M16").Value = 1.032779339687592
$ws.Range("N16").Value = 1.013142475666472

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021132602535216
$ws.Range("D17").Value = 1.027625068249576
$ws.Range("E17").Value = 1.022121267627404
$ws.Range("F17").Value = 1.030167492744716
$ws.Range("I17").Value = 1.032312358733342
$ws.Range("J17").Value = 1.027415750576016
$ws.Range("K17").Value = 1.031020357561542
$ws.Range("L17").Value = 1.025536200004149
$ws.Range("M17").Value = 1.033553786276729
$ws.Range("N17").Value = 1.013194820499157

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021305509131796
$ws.Range("D18").Value = 1.027757514385733
$ws.Range("E18").Value = 1.022266334924548
$ws.Range("F18").Value = 1.030662391033619
$ws.Range("I18").Value = 1.032361021269465
$ws.Range("J18").Value = 1.027507147135765
$ws.Range("K18").Value = 1.031110044532197
$ws.Range("L18").Value = 1.025638197769829
$ws.Range("M18").Value = 1.034004785783265
$ws.Range("N18").Value = 1.013225339786986

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021364474763595
$ws.Range("D19").Value = 1.027802680131156
$ws.Range("E19").Value = 1.022315814839384
$ws.Range("F19").Value = 1.030831022793813
$ws.Range("I19").Value = 1.032377587658852
$ws.Range("J19").Value = 1.027538306189331
$ws.Range("K19").Value = 1.031140616275563
$ws.Range("L19").Value = 1.02567297811291
$ws.Range("M19").Value = 1.034158443493007
$ws.Range("N19").Value = 1.013235743940819

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021100801920958
$ws.Range("D20").Value = 1.02760070817259
$ws.Range("E20").Value = 1.02209459105565
$ws.Range("F20").Value = 1.030076404722919
$ws.Range("I20").Value = 1.032303395144851
$ws.Range("J20").Value = 1.027398936583913
$ws.Range("K20").Value = 1.031003855952576
$ws.Range("L20").Value = 1.025517439112019
$ws.Range("M20").Value = 1.033470770418039
$ws.Range("N20").Value = 1.013189205693265

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020244199282537
$ws.Range("D21").Value = 1.02694442866318
$ws.Range("E21").Value = 1.021376464362918
$ws.Range("F21").Value = 1.027614889796046
$ws.Range("I21").Value = 1.032060362043693
$ws.Range("J21").Value = 1.026945499818894
$ws.Range("K21").Value = 1.030558597271763
$ws.Range("L21").Value = 1.025011889369392
$ws.Range("M21").Value = 1.031226516057921
$ws.Range("N21").Value = 1.013037757847645

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019706100464496
$ws.Range("D22").Value = 1.026532070779291
$ws.Range("E22").Value = 1.020925792516652
$ws.Range("F22").Value = 1.026060781035137
$ws.Range("I22").Value = 1.03190615152541
$ws.Range("J22").Value = 1.026660151847969
$ws.Range("K22").Value = 1.030278156088064
$ws.Range("L22").Value = 1.02469412638781
$ws.Range("M22").Value = 1.029808726253844
$ws.Range("N22").Value = 1.01294242397242

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01999131098646
$ws.Range("D23").Value = 1.026750643664262
$ws.Range("E23").Value = 1.021164621974677
$ws.Range("F23").Value = 1.026885267751607
$ws.Range("I23").Value = 1.031988035118495
$ws.Range("J23").Value = 1.026811444264448
$ws.Range("K23").Value = 1.030426869697029
$ws.Range("L23").Value = 1.024862569295241
$ws.Range("M23").Value = 1.030560973534706
$ws.Range("N23").Value = 1.012992972908403

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02111517108672
$ws.Range("D24").Value = 1.027611715351182
$ws.Range("E24").Value = 1.022106644763718
$ws.Range("F24").Value = 1.030117565606555
$ws.Range("I24").Value = 1.0323074458841
$ws.Range("J24").Value = 1.027406534191002
$ws.Range("K24").Value = 1.031011312487173
$ws.Range("L24").Value = 1.02552591631994
$ws.Range("M24").Value = 1.033508283941189
$ws.Range("N24").Value = 1.013191742821282

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02242142047527
$ws.Range("D25").Value = 1.028612107046173
$ws.Range("E25").Value = 1.023203447528346
$ws.Range("F25").Value = 1.033841565710865
$ws.Range("I25").Value = 1.032672022587533
$ws.Range("J25").Value = 1.028095990569142
$ws.Range("K25").Value = 1.031687398510026
$ws.Range("L25").Value = 1.026296097633164
$ws.Range("M25").Value = 1.036900262543334
$ws.Range("N25").Value = 1.013421910829964

